$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6415.909
$ws.Range("I116").Value = 5084.3
$ws.Range("J116").Value = 7525.5835
$ws.Range("K116").Value = 5084.3
$ws.Range("L116").Value = 7525.5835
$ws.Range("M116").Value = -1642.3
$ws.Range("N116").Value = -14409.5835

$ws.Range("H118").Value = 83333710
$ws.Range("I118").Value = 100000344
$ws.Range("K118").Value = 300001032
$ws.Range("M118").Value = -299999375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11090.542
$ws.Range("I32").Value = 6658.4863
$ws.Range("K32").Value = 6658.4863
$ws.Range("M32").Value = -6371.4863

$ws.Range("H61").Value = 4113.9546
$ws.Range("I61").Value = 4151.222
$ws.Range("J61").Value = 3946.25
$ws.Range("K61").Value = 4151.222
$ws.Range("L61").Value = 3946.25
$ws.Range("M61").Value = -3939.222
$ws.Range("N61").Value = -4370.25

$ws.Range("H63").Value = 5796.357
$ws.Range("I63").Value = 2625
$ws.Range("J63").Value = 8174.875
$ws.Range("K63").Value = 2625
$ws.Range("L63").Value = 8174.875
$ws.Range("M63").Value = -1939
$ws.Range("N63").Value = -9546.875

$ws.Range("H66").Value = 5796.357
$ws.Range("I66").Value = 2625
$ws.Range("J66").Value = 8174.875
$ws.Range("K66").Value = 13125
$ws.Range("L66").Value = 40874.375
$ws.Range("M66").Value = -9693
$ws.Range("N66").Value = -47738.375

$ws.Range("H136").Value = 4113.9546
$ws.Range("I136").Value = 4151.222
$ws.Range("J136").Value = 3946.25
$ws.Range("K136").Value = 12453.666
$ws.Range("L136").Value = 11838.75
$ws.Range("M136").Value = -9903.665999999999
$ws.Range("N136").Value = -16938.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 134779
$ws.Range("J57").Value = 134779
$ws.Range("L57").Value = 134779
$ws.Range("N57").Value = -136219

$ws.Range("H58").Value = 15779
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 15779
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 15779
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -16367

$ws.Range("H59").Value = 123000
$ws.Range("J59").Value = 123000
$ws.Range("L59").Value = 123000
$ws.Range("N59").Value = -124694

$ws.Range("H60").Value = 47593.332
$ws.Range("J60").Value = 47593.332
$ws.Range("L60").Value = 47593.332
$ws.Range("N60").Value = -48791.332

$ws.Range("H74").Value = 34823.8
$ws.Range("J74").Value = 34823.8
$ws.Range("L74").Value = 34823.8
$ws.Range("N74").Value = -36695.8

$ws.Range("H77").Value = 34823.8
$ws.Range("J77").Value = 34823.8
$ws.Range("L77").Value = 104471.4
$ws.Range("N77").Value = -113831.4

$ws.Range("H81").Value = 11407.25
$ws.Range("J81").Value = 11407.25
$ws.Range("L81").Value = 11407.25
$ws.Range("N81").Value = -13529.25

$ws.Range("H84").Value = 11407.25
$ws.Range("J84").Value = 11407.25
$ws.Range("L84").Value = 34221.75
$ws.Range("N84").Value = -44829.75

$ws.Range("H134").Value = 9177.111000000001
$ws.Range("I134").Value = 3148.5
$ws.Range("J134").Value = 14000
$ws.Range("K134").Value = 9445.5
$ws.Range("L134").Value = 42000
$ws.Range("M134").Value = -6910.5
$ws.Range("N134").Value = -47070

$ws.Range("H136").Value = 134779
$ws.Range("J136").Value = 134779
$ws.Range("L136").Value = 134779
$ws.Range("N136").Value = -144979

$ws.Range("H137").Value = 86854.5
$ws.Range("I137").Value = 85709
$ws.Range("J137").Value = 88000
$ws.Range("K137").Value = 85709
$ws.Range("L137").Value = 88000
$ws.Range("M137").Value = -80609
$ws.Range("N137").Value = -98200

$ws.Range("H138").Value = 65563.625
$ws.Range("J138").Value = 65563.625
$ws.Range("L138").Value = 65563.625
$ws.Range("N138").Value = -75843.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 159955
$ws.Range("I132").Value = 144512.42
$ws.Range("K132").Value = 433537.26
$ws.Range("M132").Value = -431007.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1884.3
$ws.Range("J132").Value = 1893.125
$ws.Range("L132").Value = 17038.125
$ws.Range("N132").Value = -22098.125

$ws.Range("H139").Value = 3165
$ws.Range("I139").Value = 1495
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 4485
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = 655
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3005000
$ws.Range("I11").Value = 4502500
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 4502500
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -4502361
$ws.Range("N11").Value = -10278

$ws.Range("H29").Value = 13714.143
$ws.Range("I29").Value = 9000
$ws.Range("J29").Value = 19999.666
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 19999.666
$ws.Range("M29").Value = -8710
$ws.Range("N29").Value = -20579.666

$ws.Range("H80").Value = 43702548
$ws.Range("I80").Value = 65551540
$ws.Range("J80").Value = 4569.25
$ws.Range("K80").Value = 65551540
$ws.Range("L80").Value = 4569.25
$ws.Range("M80").Value = -65550542
$ws.Range("N80").Value = -6565.25

$ws.Range("H83").Value = 43702548
$ws.Range("I83").Value = 65551540
$ws.Range("J83").Value = 4569.25
$ws.Range("K83").Value = 327757700
$ws.Range("L83").Value = 22846.25
$ws.Range("M83").Value = -327752708
$ws.Range("N83").Value = -32830.25

$ws.Range("H102").Value = 6427483
$ws.Range("I102").Value = 11112263
$ws.Range("K102").Value = 11112263
$ws.Range("M102").Value = -11110641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2643.4
$ws.Range("I55").Value = 1855.75
$ws.Range("J55").Value = 3168.5
$ws.Range("K55").Value = 1855.75
$ws.Range("L55").Value = 3168.5
$ws.Range("M55").Value = -1682.75
$ws.Range("N55").Value = -3514.5

$ws.Range("H93").Value = 37059670
$ws.Range("I93").Value = 83333736
$ws.Range("J93").Value = 40420.2
$ws.Range("K93").Value = 83333736
$ws.Range("L93").Value = 40420.2
$ws.Range("M93").Value = -83332488
$ws.Range("N93").Value = -42916.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 18000
$ws.Range("J18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("N18").Value = -18346

$ws.Range("H81").Value = 11910301
$ws.Range("I81").Value = 18521024
$ws.Range("K81").Value = 37042048
$ws.Range("M81").Value = -37040987

$ws.Range("H84").Value = 11910301
$ws.Range("I84").Value = 18521024
$ws.Range("K84").Value = 185210240
$ws.Range("M84").Value = -185204936

$ws.Range("H132").Value = 29737010
$ws.Range("I132").Value = 45461256
$ws.Range("J132").Value = 909226.0600000001
$ws.Range("K132").Value = 136383768
$ws.Range("L132").Value = 2727678.18
$ws.Range("M132").Value = -136381238
$ws.Range("N132").Value = -2732738.18
